# Repull data, push all data, mean calculation
# Update the dSF (F) column values for the skubal_tarik game log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "F2"  = -3
    "F4"  = 3
    "F5"  = 2
    "F6"  = 3
    "F7"  = 1
    "F9"  = 1
    "F10" = 6
    "F11" = -1
    "F12" = 10
    "F14" = -1
    "F15" = 1
    "F16" = 5
    "F17" = 1
    "F18" = 1
    "F20" = 4
    "F21" = -5
    "F22" = -4
    "F23" = 10
    "F25" = 8
    "F26" = -4
    "F27" = 13
    "F28" = 6
    "F29" = -2
    "F30" = 3
    "F31" = 6
    "F32" = -2
    "F33" = 2
    "F34" = 1
    "F35" = 1
    "F36" = -1
    "F37" = -2
    "F38" = 1
}

foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}
